# Rename the worksheet's tab (Sheet1 -> Лист1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Лист1"

# Insert a header row (row 1) with column titles, without shifting
# the existing data rows (which stay on rows 2 and 3).
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Full Name"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Security Question"
$ws.Range("G1").Value = "Security Answer"

# Widen the Email / Security Question / Security Answer columns to fit
# their header text (values chosen so the persisted width matches the
# target workbook's column widths of ~19.14 / ~18.86 / ~17.71 chars).
$ws.Range("E:E").ColumnWidth = 18.307291666666668
$ws.Range("F:F").ColumnWidth = 18.022135416666668
$ws.Range("G:G").ColumnWidth = 16.877604166666668

# Select the data range, leaving the header row unselected, matching
# the saved selection/active cell in the workbook.
$null = $ws.Range("A2:G3").Select()
